$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.278.84'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.11%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.865.95'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.37%  '

$ws.Range("E4").Value = '  +0.64%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.41%  '

$ws.Range("E6").Value = '  +0.83%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.49'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.81%  '

$ws.Range("E9").Value = '  +0.82%  '

$ws.Range("E10").Value = '  +1.35%  '

$ws.Range("E11").Value = '  +0.79%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.134.52'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.28%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.55'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.872.48'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.23%  '

$ws.Range("E15").Value = '  +1.24%  '

$ws.Range("E16").Value = '  +1.85%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '35.270.12'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.21'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.51%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0797'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '241.31'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.35%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.58%  '

$ws.Range("E22").Value = '  +1.47%  '

$ws.Range("E23").Value = '  +0.63%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.27'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '169.84'
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +25.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.18'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.07%  '

$ws.Range("E28").Value = '  +1.82%  '

$ws.Range("E29").Value = '  +0.88%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0565'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.12%  '

$ws.Range("E31").Value = '  +0.59%  '

$ws.Range("E32").Value = '  +2.23%  '

$ws.Range("E33").Value = '  +28.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.05'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.31%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +9.11%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.820'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +17.88%  '

$ws.Range("E37").Value = '  +6.04%  '

$ws.Range("E38").Value = '  +3.90%  '

$ws.Range("E39").Value = '  +4.45%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '90.47'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.65%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.346.39'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '15.32'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.55%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0605'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +15.83%  '

$ws.Range("E44").Value = '  +2.55%  '

$ws.Range("E45").Value = '  +0.66%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.41'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +46.87%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.64'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.45%  '

$ws.Range("E48").Value = '  -0.80%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.050.44'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.46%  '

$ws.Range("E50").Value = '  +3.48%  '

$ws.Range("E51").Value = '  +1.10%  '
